{"js": "// Remove the \"no carga autom\u00e1ticamente\" paragraph and the empty paragraph\n// that used to separate it from the \"INTERFAZ PUBLICA\" heading, since the\n// word-loading service is now synchronous and loads the word list\n// automatically during construction.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"El Servicio no carga de forma autom\";\nconst toDelete = [];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.indexOf(targetText) !== -1) {\n    // Delete the paragraph with the obsolete explanation, and the blank\n    // paragraph immediately preceding it (which separated it from the\n    // \"INTERFAZ PUBLICA\" heading).\n    toDelete.push(paragraphs.items[i]);\n    if (i > 0 && paragraphs.items[i - 1].text === \"\") {\n      toDelete.push(paragraphs.items[i - 1]);\n    }\n    break;\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"no carga autom\u00e1ticamente\" paragraph and the empty paragraph\n# that used to separate it from the \"INTERFAZ PUBLICA\" heading, since the\n# word-loading service is now synchronous and loads the word list\n# automatically during construction.\n$d = $word.ActiveDocument\n\n$targetText = \"El Servicio no carga de forma autom\"\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"$targetText*\") {\n        # Delete the paragraph with the obsolete explanation.\n        $p.Range.Delete()\n        # Delete the blank paragraph immediately before it (which separated\n        # it from the \"INTERFAZ PUBLICA\" heading), if present.\n        if ($i - 1 -ge 1) {\n            $prev = $d.Paragraphs.Item($i - 1)\n            if ($prev.Range.Text.Trim().Length -eq 0) {\n                $prev.Range.Delete()\n            }\n        }\n        break\n    }\n}\n"}
